# Auto-generated edit script applying numeric corrections to H..N columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1077.122
$ws.Range("J17").Value = 1082.8975
$ws.Range("L17").Value = 3248.6925
$ws.Range("N17").Value = -3584.6925

# Row 40
$ws.Range("H40").Value = 3197.8438
$ws.Range("I40").Value = 4987.846
$ws.Range("J40").Value = 1973.1052
$ws.Range("K40").Value = 4987.846
$ws.Range("L40").Value = 1973.1052
$ws.Range("M40").Value = -4812.846
$ws.Range("N40").Value = -2323.1052

# Row 74
$ws.Range("H74").Value = 7254.4546
$ws.Range("I74").Value = 5630.5
$ws.Range("K74").Value = 5630.5
$ws.Range("M74").Value = -4694.5

# Row 77
$ws.Range("H77").Value = 7254.4546
$ws.Range("I77").Value = 5630.5
$ws.Range("K77").Value = 28152.5
$ws.Range("M77").Value = -23472.5

# Row 100
$ws.Range("H100").Value = 5654.8945
$ws.Range("I100").Value = 2269.3333
$ws.Range("K100").Value = 2269.3333
$ws.Range("M100").Value = -1728.3333

# Row 112
$ws.Range("H112").Value = 6197.6274
$ws.Range("I112").Value = 1060
$ws.Range("J112").Value = 6756.0654
$ws.Range("K112").Value = 3180
$ws.Range("L112").Value = 20268.1962
$ws.Range("M112").Value = -2072
$ws.Range("N112").Value = -22484.1962

# Row 113
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 454
$ws.Range("N113").Value = -8708

# Row 132
$ws.Range("H132").Value = 91927.84
$ws.Range("I132").Value = 98963.5
$ws.Range("K132").Value = 296890.5
$ws.Range("M132").Value = -294360.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 23372.125
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 32
$ws.Range("H32").Value = 20962.418
$ws.Range("I32").Value = 21658.17
$ws.Range("K32").Value = 21658.17
$ws.Range("M32").Value = -21371.17

# Row 110
$ws.Range("H110").Value = 1079.5
$ws.Range("I110").Value = 427.4
$ws.Range("J110").Value = 2166.3333
$ws.Range("K110").Value = 427.4
$ws.Range("L110").Value = 2166.3333
$ws.Range("M110").Value = 1617.6
$ws.Range("N110").Value = -6256.3333

# Row 116
$ws.Range("H116").Value = 23372.125
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 23372.125
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# Row 20
$ws.Range("H20").Value = 1014.46155
$ws.Range("I20").Value = 1020
$ws.Range("J20").Value = 996
$ws.Range("K20").Value = 1020
$ws.Range("L20").Value = 996
$ws.Range("M20").Value = -773
$ws.Range("N20").Value = -1490

# Row 22
$ws.Range("H22").Value = 532.8889
$ws.Range("I22").Value = 532.8889
$ws.Range("K22").Value = 532.8889
$ws.Range("M22").Value = -359.8889

# Row 128
$ws.Range("H128").Value = 16333
$ws.Range("I128").Value = 16333
$ws.Range("K128").Value = 48999
$ws.Range("M128").Value = -46509

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 449.5
$ws.Range("I16").Value = 449
$ws.Range("K16").Value = 449
$ws.Range("M16").Value = -162

# Row 31
$ws.Range("H31").Value = 33100.176
$ws.Range("I31").Value = 23446.5
$ws.Range("J31").Value = 36070.54
$ws.Range("K31").Value = 23446.5
$ws.Range("L31").Value = 36070.54
$ws.Range("M31").Value = -23151.5
$ws.Range("N31").Value = -36660.54

# Row 34
$ws.Range("H34").Value = 33100.176
$ws.Range("I34").Value = 23446.5
$ws.Range("J34").Value = 36070.54
$ws.Range("K34").Value = 23446.5
$ws.Range("L34").Value = 36070.54
$ws.Range("M34").Value = -23244.5
$ws.Range("N34").Value = -36474.54

# Row 69
$ws.Range("H69").Value = 50000
$ws.Range("I69").Value = 12250
$ws.Range("J69").Value = 201000
$ws.Range("K69").Value = 12250
$ws.Range("L69").Value = 201000
$ws.Range("M69").Value = -11501
$ws.Range("N69").Value = -202498

# Row 72
$ws.Range("H72").Value = 50000
$ws.Range("I72").Value = 12250
$ws.Range("J72").Value = 201000
$ws.Range("K72").Value = 36750
$ws.Range("L72").Value = 603000
$ws.Range("M72").Value = -33006
$ws.Range("N72").Value = -610488

# Row 107
$ws.Range("H107").Value = 746.45715
$ws.Range("I107").Value = 679.2632
$ws.Range("J107").Value = 826.25
$ws.Range("K107").Value = 679.2632
$ws.Range("L107").Value = 826.25
$ws.Range("M107").Value = 1240.7368
$ws.Range("N107").Value = -4666.25

# Row 113
$ws.Range("H113").Value = 449.5
$ws.Range("I113").Value = 449
$ws.Range("K113").Value = 449
$ws.Range("M113").Value = 1721

# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").Value = 0

# Row 122
$ws.Range("H122").Value = 4143
$ws.Range("J122").Value = 4978.8335
$ws.Range("L122").Value = 14936.5005
$ws.Range("N122").Value = -19836.5005

# Row 132
$ws.Range("H132").Value = 240322
$ws.Range("I132").Value = 2100.6667
$ws.Range("K132").Value = 6302.000100000001
$ws.Range("M132").Value = -3772.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 900
$ws.Range("J36").Value = 1050
$ws.Range("L36").Value = 3150
$ws.Range("N36").Value = -3488

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 13251
$ws.Range("I80").Value = 21250.572
$ws.Range("J80").Value = 3918.1667
$ws.Range("K80").Value = 21250.572
$ws.Range("L80").Value = 3918.1667
$ws.Range("M80").Value = -20252.572
$ws.Range("N80").Value = -5914.1667

# Row 83
$ws.Range("H83").Value = 13251
$ws.Range("I83").Value = 21250.572
$ws.Range("J83").Value = 3918.1667
$ws.Range("K83").Value = 106252.86
$ws.Range("L83").Value = 19590.8335
$ws.Range("M83").Value = -101260.86
$ws.Range("N83").Value = -29574.8335

# Row 113
$ws.Range("H113").Value = 5239.515
$ws.Range("I113").Value = 4258.1577
$ws.Range("K113").Value = 4258.1577
$ws.Range("M113").Value = -2088.1577

# Row 132
$ws.Range("H132").Value = 10116.333
$ws.Range("J132").Value = 10866
$ws.Range("L132").Value = 32598
$ws.Range("N132").Value = -37658

$ws = $wb.Worksheets.Item("LTW")
# Row 120
$ws.Range("H120").Value = 100698
$ws.Range("J120").Value = 100698
$ws.Range("L120").Value = 100698
$ws.Range("N120").Value = -110374

# Row 122
$ws.Range("H122").Value = 4574.75
$ws.Range("I122").Value = 4205.478
$ws.Range("K122").Value = 12616.434
$ws.Range("M122").Value = -10166.434

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 752.2727
$ws.Range("I100").Value = 804
$ws.Range("K100").Value = 1608
$ws.Range("M100").Value = -1067

# Row 113
$ws.Range("H113").Value = 1884.1923
$ws.Range("I113").Value = 652.3333
$ws.Range("J113").Value = 3564
$ws.Range("K113").Value = 1956.9999
$ws.Range("L113").Value = 10692
$ws.Range("M113").Value = 213.0001
$ws.Range("N113").Value = -15032

# Row 122
$ws.Range("H122").Value = 1995.25
$ws.Range("I122").Value = 1664.2354
$ws.Range("K122").Value = 4992.706200000001
$ws.Range("M122").Value = -2542.706200000001

# Row 126
$ws.Range("H126").Value = 1758.9
$ws.Range("I126").Value = 1671.2778
$ws.Range("J126").Value = 2547.5
$ws.Range("K126").Value = 5013.8334
$ws.Range("L126").Value = 7642.5
$ws.Range("M126").Value = -2543.8334
$ws.Range("N126").Value = -12582.5

